$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 19040
$ws.Range("B2").Value = "19200_snare_04"
$ws.Range("C2").Value = "근육사슬"
$ws.Range("D2").Value = "{(snare, 2)}"
$ws.Range("E2").Value = "target"

$ws.Range("A3:E4").Clear()

$ws.Range("G7").Select()
